$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 277 (pushing existing rows 277-373 down to 279-375)
$ws.Rows("277:278").Insert()

# --- New row 277 ---
$ws.Range("A277").Value = 5
$ws.Range("B277").Value = "Macroferia Regional de Talca"
$ws.Range("C277").Value = "Maule"
$ws.Range("D277").Value = 44524
$ws.Range("E277").Value = 7
$ws.Range("F277").Value = "Fruta"
$ws.Range("G277").Value = 100101
$ws.Range("H277").Value = "Berries"
$ws.Range("I277").Value = 100112025
$ws.Range("J277").Value = "Frutilla"
$ws.Range("K277").Value = "Sin especificar"
$ws.Range("L277").Value = "Especial"
$ws.Range("M277").Value = 350
$ws.Range("N277").Value = 8000
$ws.Range("O277").Value = 8000
$ws.Range("P277").Value = 8000
$ws.Range("Q277").Value = "$/bandeja 7 kilos"
$ws.Range("R277").Value = "Provincia de Melipilla"
$ws.Range("S277").Value = 1143
$ws.Range("T277").Value = 7

# --- New row 278 ---
$ws.Range("A278").Value = 5
$ws.Range("B278").Value = "Macroferia Regional de Talca"
$ws.Range("C278").Value = "Maule"
$ws.Range("D278").Value = 44524
$ws.Range("E278").Value = 7
$ws.Range("F278").Value = "Fruta"
$ws.Range("G278").Value = 100101
$ws.Range("H278").Value = "Berries"
$ws.Range("I278").Value = 100112025
$ws.Range("J278").Value = "Frutilla"
$ws.Range("K278").Value = "Sin especificar"
$ws.Range("L278").Value = "Segunda"
$ws.Range("M278").Value = 150
$ws.Range("N278").Value = 5000
$ws.Range("O278").Value = 5000
$ws.Range("P278").Value = 5000
$ws.Range("Q278").Value = "$/bandeja 7 kilos"
$ws.Range("R278").Value = "Provincia de Melipilla"
$ws.Range("S278").Value = 714
$ws.Range("T278").Value = 7
